$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Update the Approved/Rejected column (I) to "Approved" and clear the
# ReasonToReject column (J) for every row that was previously "Rejected"/"Nil".
foreach ($r in 2,4,5,6,7,8,9) {
    $ws.Range("I$r").Value = "Approved"
    $ws.Range("J$r").ClearContents() | Out-Null
}

# Update the active selection shown in the sheet view.
$ws.Range("I17").Select() | Out-Null
